# Apply the PNAD 2009 "roubo/furto" sheet fixes:
#  1. Rename the header cell B2 from "unnamed: 1_level_1" to "total".
#  2. Remove the two empty section-header rows ("situação do domicílio" and
#     "grandes regiões e unidades da federação"), which had no data of
#     their own - deleting them shifts every row below up so the real data
#     lines up with its label again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the mis-labelled sub-header.
$ws.Range("B2").Value = "total"

# 2) Delete row 5 ("situação do domicílio") - a blank separator row.
$ws.Rows(5).Delete()

# 3) After the row-5 delete, the other blank separator row
#    ("grandes regiões e unidades da federação") has shifted from row 8
#    up to row 7 - delete it too.
$ws.Rows(7).Delete()
